# Add a new "2021" data column (column R) to the right of the existing
# "2020" column (column Q), mirroring the formatting of column Q for each
# row, then set the new value for each row and leave the selection on R2
# (matching the author's final cursor position after entering the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (thin separator row): empty cell, same style as Q2.
$ws.Range("Q2").Copy($ws.Range("R2"))

# Row 3 (year header row): 2021 follows 2020.
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("R3").Value = 2021

# Row 4 (population count row).
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 202551

# Row 5 (percentage row).
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 2.9794303052841493

# Leave the selection where the author left it after editing the sheet.
$ws.Range("R2").Select() | Out-Null
